$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 2 to hold the data for years 2018, 2019 and 2020.
# This shifts the existing 2021-2024 rows down from rows 2-5 to rows 5-8.
$ws.Rows.Item(2).Resize(3).Insert()

# Match the formatting (bold, centered, bordered, text format) already used
# for the year column (A) in the other rows.
$yearRange = $ws.Range("A2:A4")
$yearRange.NumberFormat = "@"
$yearRange.Font.Bold = $true
$yearRange.HorizontalAlignment = -4108  # xlCenter
$yearRange.VerticalAlignment = -4160    # xlTop
$yearRange.Borders.LineStyle = 1        # xlContinuous

# Row 2 -> 2018
$ws.Range("A2").Value = "2018"
$ws.Range("B2").Value = 0.05
$ws.Range("C2").Value = 24.32
$ws.Range("D2").Value = 6.84
$ws.Range("E2").Value = 60.98
$ws.Range("F2").Value = 45.64
$ws.Range("G2").Value = 99.52
$ws.Range("H2").Value = 90.01000000000001
$ws.Range("I2").Value = 85.3
$ws.Range("J2").Value = 37.73
$ws.Range("K2").Value = 0.59
$ws.Range("L2").Value = 0.91
$ws.Range("M2").Value = 42.81
$ws.Range("N2").Value = 0.91
$ws.Range("O2").Value = 70.81999999999999
$ws.Range("P2").Value = 99.2

# Row 3 -> 2019
$ws.Range("A3").Value = "2019"
$ws.Range("B3").Value = 0.5600000000000001
$ws.Range("C3").Value = 32.22
$ws.Range("D3").Value = 4.37
$ws.Range("E3").Value = 98.09999999999999
$ws.Range("F3").Value = 31.51
$ws.Range("G3").Value = 99.05
$ws.Range("H3").Value = 83.33
$ws.Range("I3").Value = 62.54
$ws.Range("J3").Value = 36.9
$ws.Range("K3").Value = 9.84
$ws.Range("L3").Value = 6.03
$ws.Range("M3").Value = 48.89
$ws.Range("N3").Value = 5.4
$ws.Range("O3").Value = 64.84
$ws.Range("P3").Value = 98.97

# Row 4 -> 2020
$ws.Range("A4").Value = "2020"
$ws.Range("B4").Value = 0.44
$ws.Range("C4").Value = 49.85
$ws.Range("D4").Value = 4.69
$ws.Range("E4").Value = 99.27
$ws.Range("F4").Value = 45.16
$ws.Range("G4").Value = 99.70999999999999
$ws.Range("H4").Value = 80.65000000000001
$ws.Range("I4").Value = 64.37
$ws.Range("J4").Value = 37.1
$ws.Range("K4").Value = 12.61
$ws.Range("L4").Value = 12.32
$ws.Range("M4").Value = 59.68
$ws.Range("N4").Value = 6.16
$ws.Range("O4").Value = 70.67
$ws.Range("P4").Value = 99.56

# The previously existing rows (now at 5-8, for years 2021-2024) only changed
# their column C ("A") value.
$ws.Range("C5").Value = 61.02
$ws.Range("C6").Value = 62.29
$ws.Range("C7").Value = 55.56
$ws.Range("C8").Value = 43.8

Write-Output "applied criterios_por_ano update"
